$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price/volume data to match the latest snapshot.
# Columns B/C are only touched where two coins swapped ranking rows (30/32, 41/42).
# Price (D) cells must stay text (e.g. "34.411.58", "1.00", "0.0500") rather than
# being auto-converted to numbers by Excel, so NumberFormat is forced to Text first.

$ws.Range("D2").Value = "34.411.58"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.801.81"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.38"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.11"
$ws.Range("E8").Value = "  +4.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.291"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0678"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").Value = "2.060.59"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.23"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "1.802.38"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.629"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").Value = "34.371.12"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.70"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.05"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "0.0₃0775"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +5.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.11"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.85"
$ws.Range("E26").Value = "  +5.02%  "
$ws.Range("E27").Value = "  +4.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.121"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.91"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.22"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0515"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("D35").Value = "1.363.15"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.649"
$ws.Range("E36").Value = "  -3.22%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -7.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0187"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "81.24"
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.79"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.936"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  +5.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.21"
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0500"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").Value = "1.963.67"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.78"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("E51").Value = "  -1.13%  "
